# Updates the cryptos list sheet with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "195.35") must be
# forced to Text format first, otherwise Excel would silently convert them
# into numeric values (dropping formatting such as trailing zeros).
$textForcedUpdates = [ordered]@{
    "D5" = '195.35'
    "D6" = '597.40'
    "D7" = '1.00'
    "D11" = '0.397'
    "D17" = '27.24'
    "D19" = '8.92'
    "D20" = '12.55'
    "D21" = '380.94'
    "D24" = '71.59'
    "D27" = '4.21'
    "D28" = '9.71'
    "D32" = '506.78'
    "D33" = '7.76'
    "D35" = '1.00'
    "D36" = '164.54'
    "D37" = '20.17'
    "D38" = '19.67'
    "D40" = '182.80'
    "D47" = '40.45'
    "D49" = '0.579'
    "D50" = '0.668'
    "D51" = '3.76'
}

foreach ($cellRef in $textForcedUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$cellRef]
}

# Remaining cells (coin names, links, already-textual prices, percentages)
# can be written directly.
$plainUpdates = [ordered]@{
    "D2" = '75.813.93'
    "E2" = '  +0.09%  '
    "D3" = '2.878.58'
    "E3" = '  +5.61%  '
    "E4" = '  +0.10%  '
    "E5" = '  +3.12%  '
    "E6" = '  +0.83%  '
    "E7" = '  +0.13%  '
    "E8" = '  +1.69%  '
    "E9" = '  -3.49%  '
    "D10" = '2.880.82'
    "E10" = '  +5.75%  '
    "E11" = '  +9.26%  '
    "E12" = '  -1.53%  '
    "E13" = '  +2.23%  '
    "D14" = '3.412.19'
    "E14" = '  +7.19%  '
    "D15" = '75.753.53'
    "E15" = '  +0.25%  '
    "E16" = '  -1.10%  '
    "E17" = '  +0.94%  '
    "D18" = '2.898.74'
    "E18" = '  +7.44%  '
    "E19" = '  -7.13%  '
    "E20" = '  +2.48%  '
    "E21" = '  +0.11%  '
    "E22" = '  -1.11%  '
    "E23" = '  +0.52%  '
    "E24" = '  +0.76%  '
    "E25" = '  -0.04%  '
    "D26" = '3.042.64'
    "E27" = '  -1.27%  '
    "E28" = '  +0.52%  '
    "E29" = '  +8.54%  '
    "E30" = '  -0.36%  '
    "E31" = '  -1.96%  '
    "E32" = '  -4.35%  '
    "E33" = '  -2.22%  '
    "E34" = '  +1.64%  '
    "E35" = '  +0.18%  '
    "E36" = '  +1.44%  '
    "E37" = '  +2.98%  '
    "E38" = '  +1.46%  '
    "E39" = '  -6.04%  '
    "E40" = '  +4.71%  '
    "E41" = '  -0.09%  '
    "E42" = '  +2.38%  '
    "E43" = '  -2.40%  '
    "E44" = '  -3.03%  '
    "E45" = '  +7.01%  '
    "E46" = '  -0.60%  '
    "E47" = '  +3.08%  '
    "E48" = '  -3.99%  '
    "E49" = '  +5.18%  '
    "B50" = 'Mantle'
    "C50" = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    "E50" = '  +11.67%  '
    "B51" = 'Filecoin'
    "C51" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    "E51" = '  +1.32%  '
}

foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}
